# Regen save_data to use K instead of Strike#, write new s_vals into column G ("K")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K (strikeout) value, per regenerated save_data
$kValues = @{
    2 = 1
    3 = 1
    4 = 0
    5 = 2
    6 = 0
    7 = 1
    8 = 1
    9 = 2
    10 = 1
    11 = 1
    12 = 3
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 2
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 2
    32 = 0
    33 = 1
    34 = 1
    35 = 2
    36 = 2
    37 = 0
    38 = 2
    39 = 0
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 3
    45 = 1
    46 = 2
    47 = 1
    48 = 2
    49 = 0
    50 = 2
    51 = 1
    52 = 2
    53 = 1
    54 = 1
    55 = 1
    56 = 2
    57 = 1
    58 = 1
    59 = 3
    60 = 1
    61 = 1
    62 = 1
    64 = 1
    65 = 0
    66 = 3
    67 = 0
    68 = 1
    69 = 1
    70 = 2
    72 = 2
    74 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
